$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Capture long, special-character-laden code snippets from the pythonCode
# sheet BEFORE we mutate its layout, so we can reuse them verbatim on the
# new "Practice Qns" sheet without retyping them (avoids any transcription
# mistakes with embedded newlines / backslashes / quotes).
# ---------------------------------------------------------------------------
$wsPy = $wb.Worksheets.Item(2)
$codeSearch          = $wsPy.Range("A4").Value2
$codeFindMaxOnes     = $wsPy.Range("A6").Value2
$codeFindNumbers     = $wsPy.Range("A8").Value2
$codeSortedSquares   = $wsPy.Range("A10").Value2

# ---------------------------------------------------------------------------
# 1) Rename "PythonCode" -> "pythonCode"
# ---------------------------------------------------------------------------
$wsPy.Name = "pythonCode"

# ---------------------------------------------------------------------------
# 2) Remove the stray selection marker left on loginSheet
# ---------------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item(1)
$wsLogin.Range("A1").Select()

# ---------------------------------------------------------------------------
# 3) Restructure the pythonCode sheet: insert a new TestId column at A,
#    pushing the old data (and its widths/bestFit formatting) one column
#    to the right automatically.
# ---------------------------------------------------------------------------
$wsPy.Columns.Item(1).Insert()

$wsPy.Range("A1").Value = "TestId"
$wsPy.Range("A2").Value = "TC001"
$wsPy.Range("A3").Value = "TC002"
$wsPy.Range("A4").Value = "TC003"

# B2 used to hold the stray print("hello"); string - no longer needed, but
# keep the cell (and its formatting) blank.
$wsPy.Range("B2").ClearContents()

# B4 becomes the actual python snippet for TC003 instead of the duplicated
# "def search..." block.
$wsPy.Range("B4").Value = "print(""hello"")"
$wsPy.Range("B4").Style = "Normal"

# Rows 5-11 no longer repeat the snippet text in column B - only the
# Result column (now C) remains.
$wsPy.Range("B5:B11").Clear()

$wsPy.Columns.Item(1).ColumnWidth = 10.86

$wsPy.Range("B4").Select()

# ---------------------------------------------------------------------------
# 4) Add the new "Practice Qns" worksheet after pythonCode
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPractice = $wb.Worksheets.Add($null, $lastSheet)
$wsPractice.Name = "Practice Qns"

$wsPractice.Range("A1").Value = "TestCaseId"
$wsPractice.Range("B1").Value = "pythonCode"

$wsPractice.Range("A2").Value = "TC001"

$wsPractice.Range("A3").Value = "TC002"
$wsPractice.Range("B3").Value = "hello"

$wsPractice.Range("A4").Value = "TC003"
$wsPractice.Range("B4").Value = "print(""Hello"")"

$wsPractice.Range("B6").Value = $codeSearch
$wsPractice.Range("B7").Value = $codeSearch
$wsPractice.Range("B8").Value = $codeFindMaxOnes
$wsPractice.Range("B9").Value = $codeFindMaxOnes
$wsPractice.Range("B10").Value = $codeFindNumbers
$wsPractice.Range("B11").Value = $codeFindNumbers
$wsPractice.Range("B12").Value = $codeSortedSquares
$wsPractice.Range("B13").Value = $codeSortedSquares

$wsPractice.Range("B6:B13").Style = "Normal 2"

$wsPractice.Columns.Item(1).ColumnWidth = 9.86
$wsPractice.Columns.Item(2).ColumnWidth = 243.42578125

$wsPractice.Range("A2").Select()

$wsPy.Activate()
